# Generate Report for handoff
#
# The "92342a75-..." entry just had a new handoff generated, so its
# "Latest Handoff Datetime" cell (column D, row 4) advances on both the
# zh-cn and the de-de status sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-02-16 09:28:16"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-02-16 09:28:28"
